$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-08-18 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-19 Tuesday", 2)

# The document contains a single table of division problems. Update each
# cell's value by addressing it positionally (row, column) so that
# duplicate "before" values (e.g. "86÷4=" appears twice) are each mapped
# to their correct, distinct replacement.
$tbl = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cellRange = $table.Cell($row, $col).Range
    $cellRange.MoveEnd(1, -1) | Out-Null   # 1 = wdCharacter; trim the trailing cell-mark
    $cellRange.Text = $text
}

# Row 1
Set-CellText $tbl 1 1 "98÷3="
Set-CellText $tbl 1 2 "14÷6="
Set-CellText $tbl 1 3 "85÷5="
Set-CellText $tbl 1 4 "55÷9="
Set-CellText $tbl 1 5 "96÷4="

# Row 5
Set-CellText $tbl 5 1 "71÷4="
Set-CellText $tbl 5 2 "96÷7="
Set-CellText $tbl 5 3 "13÷4="
Set-CellText $tbl 5 4 "92÷4="
Set-CellText $tbl 5 5 "95÷8="

# Row 9
Set-CellText $tbl 9 1 "84÷6="
Set-CellText $tbl 9 2 "73÷7="
Set-CellText $tbl 9 3 "57÷3="
Set-CellText $tbl 9 4 "57÷4="
Set-CellText $tbl 9 5 "68÷6="

# Row 13
Set-CellText $tbl 13 1 "82÷9="
Set-CellText $tbl 13 2 "46÷2="
Set-CellText $tbl 13 3 "72÷8="
# Cell (13,4) "65÷2=" is unchanged.
Set-CellText $tbl 13 5 "75÷8="

# Row 17
Set-CellText $tbl 17 1 "26÷9="
Set-CellText $tbl 17 2 "69÷2="
Set-CellText $tbl 17 3 "79÷9="
Set-CellText $tbl 17 4 "73÷9="
Set-CellText $tbl 17 5 "33÷7="
